$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 2")

# Row pairs whose data (columns B through AD) needs to be swapped.
# Column A (row index number) is left untouched.
$pairs = @(
    @(13, 15),
    @(66, 67),
    @(70, 72),
    @(84, 85),
    @(118, 119),
    @(160, 161),
    @(193, 194)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$($r1):AD$($r1)")
    $range2 = $ws.Range("B$($r2):AD$($r2)")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
